# Insert the newest (3rd) tender item as a new row 2, pushing the
# existing rows down. Column D (Procurement Category) repeats the same
# value as the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: shift existing rows down, starting at row 2.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "INVITATION TO QUOTE FOR BOOK VOUCHERS"
$ws.Range("B2").Value = "19 May 2021`n01:00PM"
$ws.Range("C2").Value = "Ministry of Education"
$ws.Range("D2").Value = "Administration & Training ⇒ Gifts & Souvenirs"

# Match formatting used by the other data rows (wrapped date column,
# same row height).
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(3).RowHeight
